$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.989.22'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '''1.917.90'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").Value = '''324.81'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''0.4602'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.3827'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '''0.07705'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").Value = '''0.9808'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '''22.21'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '''1.934.74'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '''5.691'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = '''6.965'
$ws.Range("D15").Value = '''0.06983'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '''84.25'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").Value = '''0.000009462'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = '''16.68'
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").Value = '''28.965.51'
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").Value = '''5.338'
$ws.Range("E22").Value = '  -1.91%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '''2.159.14'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = '''2.092'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("D27").Value = '''19.04'
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").Value = '''5.702'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '''117.79'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '''1.866'
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").Value = '''0.09313'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '''0.8668'
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").Value = '''5.113'
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = '''1.251'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").Value = '''3.044'
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '''1.157'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '''0.02043'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").Value = '''3.034'
$ws.Range("E40").Value = '  +11.61%  '
$ws.Range("D41").Value = '''7.526'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '''0.5503'
$ws.Range("E42").Value = '  -1.31%  '
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").Value = '''9.404'
$ws.Range("E44").Value = '  +2.33%  '
$ws.Range("D45").Value = '''0.000002868'
$ws.Range("E45").Value = '  -3.90%  '
$ws.Range("D46").Value = '''2.184'
$ws.Range("E46").Value = '  +5.95%  '
$ws.Range("D47").Value = '''0.5185'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '''11.21'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").Value = '''0.06898'
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").Value = '''1.782'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").Value = '''110.40'
$ws.Range("E51").Value = '  -0.47%  '
